# ============================================================
# edit.ps1 - applies "Beginnings of documentation for Elex format"
# ============================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# --- Clear all existing content so we can rebuild the sheet from scratch ---
$ws1.Cells.ClearContents()

# --- Step 1: re-establish column A values (rows 2-41) in original order so
#     the existing shared-string indices (0-39) are preserved ---
$ws1.Range("A2").Value = 'id'
$ws1.Range("A3").Value = 'raceid'
$ws1.Range("A4").Value = 'racetype'
$ws1.Range("A5").Value = 'racetypeid'
$ws1.Range("A6").Value = 'ballotorder'
$ws1.Range("A7").Value = 'candidateid'
$ws1.Range("A8").Value = 'description'
$ws1.Range("A9").Value = 'delegatecount'
$ws1.Range("A10").Value = 'electiondate'
$ws1.Range("A11").Value = 'electtotal'
$ws1.Range("A12").Value = 'electwon'
$ws1.Range("A13").Value = 'fipscode'
$ws1.Range("A14").Value = 'first'
$ws1.Range("A15").Value = 'incumbent'
$ws1.Range("A16").Value = 'initialization_data'
$ws1.Range("A17").Value = 'is_ballot_measure'
$ws1.Range("A18").Value = 'last'
$ws1.Range("A19").Value = 'lastupdated'
$ws1.Range("A20").Value = 'level'
$ws1.Range("A21").Value = 'national'
$ws1.Range("A22").Value = 'officeid'
$ws1.Range("A23").Value = 'officename'
$ws1.Range("A24").Value = 'party'
$ws1.Range("A25").Value = 'polid'
$ws1.Range("A26").Value = 'polnum'
$ws1.Range("A27").Value = 'precinctsreporting'
$ws1.Range("A28").Value = 'precinctsreportingpct'
$ws1.Range("A29").Value = 'precinctstotal'
$ws1.Range("A30").Value = 'reportingunitid'
$ws1.Range("A31").Value = 'reportingunitname'
$ws1.Range("A32").Value = 'runoff'
$ws1.Range("A33").Value = 'seatname'
$ws1.Range("A34").Value = 'seatnum'
$ws1.Range("A35").Value = 'statename'
$ws1.Range("A36").Value = 'statepostal'
$ws1.Range("A37").Value = 'test'
$ws1.Range("A38").Value = 'uncontested'
$ws1.Range("A39").Value = 'votecount'
$ws1.Range("A40").Value = 'votepct'
$ws1.Range("A41").Value = 'winner'

# --- Step 2: add the new header row + field descriptions, in the exact order
#     they were originally authored (controls new shared-string ordering) ---
$ws1.Range("A1").Value = 'FieldName'
$ws1.Range("B1").Value = 'Description'
$ws1.Range("B3").Value = 'Unique race ID for a specific state'
$ws1.Range("B4").Value = 'Character string indicating the type of race (for example, GOP Primary, General, Democratic Caucus)'
$ws1.Range("B5").Value = 'Single-character race type ID D (Dem Primary), R (GOP Primary), G (General), E (Dem Caucus), S (GOP Caucus)'
$ws1.Range("B2").Value = '? In practice, raceid and delimiter and reportingunitid'
$ws1.Range("B7").Value = 'AP-assigned unique ID for this candidate in a state''s race. If a candidate is running in multiple races, this candidate has a different candidateid in each race'
$ws1.Range("B8").Value = 'Description of the office, ballot initiative or other (if applicable)'
$ws1.Range("B6").Value = 'Ballot order of this candidate. There may be gaps in sequence in this order field.'
$ws1.Range("B9").Value = 'For presidential primaries, delegates won by this candidate in this district'
$ws1.Range("B11").Value = 'In a general election in a presidential year, the state or U.S. national electoral count'
$ws1.Range("B12").Value = 'In a general election in a presidential year, the candidate''s electoral votes'
$ws1.Range("B13").Value = 'County FIPS code, a geographical standard that allows data to be matched'
$ws1.Range("B14").Value = 'Candidate''s first name'
$ws1.Range("B15").Value = 'A flag that indicates the candidate is an incumbent. In Elex CSV format, this is TRUE or FALSE'
$ws1.Range("B16").Value = '?'
$ws1.Range("B18").Value = 'Candidate''s last name'
$ws1.Range("B10").Value = 'Date of the election day, in format m/d/yyyy'
$ws1.Range("B19").Value = 'Time last updated. Sample: 2018-03-07T17:24:47.600Z'
$ws1.Range("B20").Value = 'Region the vote results are reported from: "national" for presidential results and electoral count at the U.S. national rolled-up summary level (only for the general election in a presidential year); "state" for state-level results; "subunit" for results at the RU or FIPS code level; "district" for delegate results at the district level from the presidential primaries or presidential results and electoral counts by district (currently, in Maine and Nebraska) for the general election in a presidential year'
$ws1.Range("B21").Value = 'Indicates that the race is national. "TRUE"'
$ws1.Range("B17").Value = '?'  # reuses '?' string from B16

# --- Formatting: bold the placeholder / TBD description cells ---
$ws1.Range("B2").Font.Bold = $true
$ws1.Range("B16").Font.Bold = $true
$ws1.Range("B17").Font.Bold = $true

# Cells that were bolded and then un-bolded by the original author
$ws1.Range("B6").Font.Bold = $true
$ws1.Range("B18:B21").Font.Bold = $true
$ws1.Range("B6").Font.Bold = $false
$ws1.Range("B18:B21").Font.Bold = $false

# --- Column widths ---
$ws1.Columns("A").ColumnWidth = 27.666666666666668
$ws1.Columns("B").ColumnWidth = 27.333333333333336

# --- Page setup: portrait orientation ---
$ws1.PageSetup.Orientation = 1

# --- Selection matches the saved view in the target workbook ---
$ws1.Range("B15").Select()

# --- Rename sheet (author renamed Sheet1 -> Sheet2) ---
$ws1.Name = "Sheet2"
